# Income statement sheet ("Overview"): prepend 5 earlier quarterly columns
# (quarters ended 1399/06 .. 1400/06) before the existing D:H block, which
# shifts the previously-reported quarters (1400/09 .. 1401/09) from D:H to
# I:M. Also refresh a handful of republished figures and their "published"
# (تاریخ انتشار) dates. Matches commit "add create_cumulative_data change
# predict_income".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Insert 5 blank columns at D, pushing the existing D:H block to I:M.
# Cell styles/number formats on the source columns are carried along.
$ws.Columns("D:H").Insert()

# --- New quarter-ended header labels (row 8) -------------------------------
$ws.Range("D8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"

# --- New "تاریخ انتشار" (publish date) row (row 9) -------------------------
$ws.Range("D9").Value = "1400-08-29 (4)"
$ws.Range("E9").Value = "1400-11-02 (2)"
$ws.Range("F9").Value = "1401-02-25 (12)"
$ws.Range("G9").Value = "1401-05-01 (3)"
$ws.Range("H9").Value = "1401-08-30 (4)"

# Republished dates for the quarters that shifted into I:M.
$ws.Range("I9").Value = "1401-11-19 (3)"
$ws.Range("J9").Value = "1401-11-19 (9)"
$ws.Range("M9").Value = "1401-11-19 (2)"

# --- New quarterly figures (rows 11-27) ------------------------------------
# فروش
$ws.Range("D11").Value = 1531119
$ws.Range("E11").Value = 2190426
$ws.Range("F11").Value = 1987544
$ws.Range("G11").Value = 2671724
$ws.Range("H11").Value = 2575491
$ws.Range("J11").Value = 3013415

# بهای تمام شده کالای فروش رفته
$ws.Range("D12").Value = -688760
$ws.Range("E12").Value = -830393
$ws.Range("F12").Value = -792083
$ws.Range("G12").Value = -1072057
$ws.Range("H12").Value = -1076741

# سود (زیان) ناخالص
$ws.Range("D13").Value = 842359
$ws.Range("E13").Value = 1360033
$ws.Range("F13").Value = 1195461
$ws.Range("G13").Value = 1599667
$ws.Range("H13").Value = 1498750
$ws.Range("J13").Value = 1285943

# هزینه های عمومی, اداری و تشکیلاتی
$ws.Range("D14").Value = -73637
$ws.Range("E14").Value = -220301
$ws.Range("F14").Value = -432882
$ws.Range("G14").Value = -290181
$ws.Range("H14").Value = -190044

# هزینه کاهش ارزش دریافتنی‌‏ها (هزینه استثنایی)
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("D16").Value = 42412
$ws.Range("E16").Value = 61880
$ws.Range("F16").Value = -25759
$ws.Range("G16").Value = -1761
$ws.Range("H16").Value = -46075

# سود (زیان) عملیاتی
$ws.Range("D17").Value = 811134
$ws.Range("E17").Value = 1201612
$ws.Range("F17").Value = 736820
$ws.Range("G17").Value = 1307725
$ws.Range("H17").Value = 1262631
$ws.Range("J17").Value = 1099263

# هزینه های مالی
$ws.Range("D18").Value = -8775
$ws.Range("E18").Value = -12804
$ws.Range("F18").Value = -9412
$ws.Range("G18").Value = -13454
$ws.Range("H18").Value = -12095

# خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("D19").Value = 50256
$ws.Range("E19").Value = 150230
$ws.Range("F19").Value = 20454
$ws.Range("G19").Value = 105926
$ws.Range("H19").Value = 65576

# سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("D20").Value = 852615
$ws.Range("E20").Value = 1339038
$ws.Range("F20").Value = 747862
$ws.Range("G20").Value = 1400197
$ws.Range("H20").Value = 1316112
$ws.Range("J20").Value = 1195388

# مالیات
$ws.Range("D21").Value = -129125
$ws.Range("E21").Value = -193764
$ws.Range("F21").Value = 21851
$ws.Range("G21").Value = -99622
$ws.Range("H21").Value = -362256

# سود (زیان) خالص عملیات در حال تداوم
$ws.Range("D22").Value = 723490
$ws.Range("E22").Value = 1145274
$ws.Range("F22").Value = 769713
$ws.Range("G22").Value = 1300575
$ws.Range("H22").Value = 953856
$ws.Range("J22").Value = 1031333

# سود (زیان) عملیات متوقف شده پس از اثر مالیاتی
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# سود (زیان) خالص
$ws.Range("D24").Value = 723490
$ws.Range("E24").Value = 1145274
$ws.Range("F24").Value = 769713
$ws.Range("G24").Value = 1300575
$ws.Range("H24").Value = 953856
$ws.Range("J24").Value = 1031333

# سود هر سهم پس از کسر مالیات
$ws.Range("D25").Value = 1113
$ws.Range("E25").Value = 1762
$ws.Range("F25").Value = 1184
$ws.Range("G25").Value = 2001
$ws.Range("H25").Value = 1467

# سرمایه
$ws.Range("D26").Value = 650000
$ws.Range("E26").Value = 650000
$ws.Range("F26").Value = 650000
$ws.Range("G26").Value = 650000
$ws.Range("H26").Value = 650000

# سود هر سهم بر اساس آخرین سرمایه
$ws.Range("D27").Value = 1113
$ws.Range("E27").Value = 1762
$ws.Range("F27").Value = 1184
$ws.Range("G27").Value = 2001
$ws.Range("H27").Value = 1467

# --- Column widths for the newly-inserted D:H block -------------------------
# (ColumnWidth is in "characters"; the saved OOXML width is ColumnWidth +
# 5/6, so subtract 5/6 from the desired OOXML width of 29/29/31/29/29.)
$ws.Range("D1").ColumnWidth = 28.16666666666667
$ws.Range("E1").ColumnWidth = 28.16666666666667
$ws.Range("F1").ColumnWidth = 30.16666666666667
$ws.Range("G1").ColumnWidth = 28.16666666666667
$ws.Range("H1").ColumnWidth = 28.16666666666667
